$wb = $excel.ActiveWorkbook

# --- Update the "outputs" sheet ---
$ws = $wb.Worksheets.Item("outputs")

# Insert two new rows before the current row 4, shifting the existing
# rows (old 4..9) down to (6..9) [and spilling the former last two rows off].
$ws.Rows.Item(4).Resize(2).Insert()

# Fill in the two newly inserted rows with the new raw_3 values.
$ws.Cells.Item(4, 1).Value = "raw_3_first_scoring_"
$ws.Cells.Item(4, 2).Value = 29
$ws.Cells.Item(5, 1).Value = "raw_3_second_scoring_"
$ws.Cells.Item(5, 2).Value = 97

# Trim the sheet back down to 9 rows (drop what is now rows 10 and 11,
# the old raw_3_first_scoring_/raw_3_second_scoring_ rows with decimal
# values), matching the target layout (dimension stays A1:C9).
$ws.Rows.Item(10).Resize(2).Delete()

# --- Add the new "outputs1" sheet, placed after "outputs" ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "outputs1"
$newSheet.Move($null, $wb.Worksheets.Item("outputs"))
